# "Generate Report for Handback" - updates the localization-status report
# after a handback has been processed for the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/e041d4e93ecc40c76a16a49b2e1063742e79ad44/e2e/a.md"

# --- Overview sheet: refresh the per-language status text -------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C3").Value = $statusHandedBack

$zhcn.Range("H2").Value = "a.md"
$zhcn.Range("H2").Style = $zhcn.Range("A2").Style
$zhcn.Hyperlinks.Add($zhcn.Range("H2"), $aMdUrl, "", "", "a.md") | Out-Null
$zhcn.Range("H2").Style = $zhcn.Range("A2").Style

$zhcn.Range("H3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("H3"), $aMdUrl, "", "", "a.md") | Out-Null
$zhcn.Range("H3").Style = $zhcn.Range("A2").Style

$zhcn.Range("I2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("I3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("J2").Value = "2016-07-26 07:42:43"
$zhcn.Range("J3").Value = "2016-07-26 07:42:43"

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664

# --- de-de sheet ----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C3").Value = $statusHandedBack

$dede.Range("H2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("H2"), $aMdUrl, "", "", "a.md") | Out-Null
$dede.Range("H2").Style = $dede.Range("A2").Style

$dede.Range("H3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("H3"), $aMdUrl, "", "", "a.md") | Out-Null
$dede.Range("H3").Style = $dede.Range("A2").Style

$dede.Range("I2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("I3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("J2").Value = "2016-07-26 07:42:59"
$dede.Range("J3").Value = "2016-07-26 07:42:59"

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
